$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking text (prices like "213.81"); Excel would
# auto-coerce plain .Value assignments of such strings into numbers, so we
# force the Text format while writing, then restore the default style.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "26.514.03"
$ws.Range("E2").Value = "  -0.81%  "
$ws.Range("D3").Value = "1.625.47"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "213.81"
$ws.Range("E5").Value = "  -0.36%  "
$ws.Range("D6").Value = "0.501"
$ws.Range("E6").Value = "  -1.07%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  -0.31%  "
$ws.Range("D9").Value = "0.0609"
$ws.Range("E9").Value = "  -0.44%  "
$ws.Range("D10").Value = "19.21"
$ws.Range("E10").Value = "  -1.05%  "
$ws.Range("E11").Value = "  -0.34%  "
$ws.Range("D12").Value = "1.855.47"
$ws.Range("E12").Value = "  +0.12%  "
$ws.Range("D13").Value = "1.614.82"
$ws.Range("E13").Value = "  -0.64%  "
$ws.Range("E14").Value = "  -0.36%  "
$ws.Range("D15").Value = "0.511"
$ws.Range("E15").Value = "  -0.43%  "
$ws.Range("D16").Value = "63.91"
$ws.Range("E16").Value = "  -1.73%  "
$ws.Range("D17").Value = "234.62"
$ws.Range("E17").Value = "  -0.05%  "
$ws.Range("D18").Value = "26.532.39"
$ws.Range("E18").Value = "  -0.79%  "
$ws.Range("E19").Value = "  +0.14%  "
$ws.Range("D20").Value = "0.0₃0725"
$ws.Range("E20").Value = "  -0.51%  "
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("E22").Value = "  -1.99%  "
$ws.Range("D23").Value = "2.21"
$ws.Range("E23").Value = "  -2.24%  "
$ws.Range("E24").Value = "  +0.28%  "
$ws.Range("D25").Value = "145.78"
$ws.Range("E25").Value = "  -0.11%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("E28").Value = "  -0.88%  "
$ws.Range("E29").Value = "  -0.49%  "
$ws.Range("E30").Value = "  -1.12%  "
$ws.Range("E31").Value = "  -0.34%  "
$ws.Range("D32").Value = "1.525.52"
$ws.Range("E32").Value = "  +3.64%  "
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("E34").Value = "  +0.06%  "
$ws.Range("D35").Value = "1.52"
$ws.Range("E35").Value = "  +2.10%  "
$ws.Range("E36").Value = "  -0.29%  "
$ws.Range("D37").Value = "0.569"
$ws.Range("E37").Value = "  -0.54%  "
$ws.Range("E38").Value = "  -0.90%  "
$ws.Range("E39").Value = "  -0.57%  "
$ws.Range("E40").Value = "  -2.07%  "
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("E42").Value = "  +0.26%  "
$ws.Range("D43").Value = "1.766.63"
$ws.Range("E43").Value = "  +0.09%  "
$ws.Range("D44").Value = "62.85"
$ws.Range("E44").Value = "  +1.19%  "
$ws.Range("D45").Value = "0.760"
$ws.Range("E45").Value = "  -1.04%  "
$ws.Range("D46").Value = "0.907"
$ws.Range("E46").Value = "  -4.55%  "
$ws.Range("D47").Value = "89.68"
$ws.Range("E47").Value = "  +1.35%  "
$ws.Range("E48").Value = "  +0.31%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₆0102"
$ws.Range("E49").Value = "  -0.97%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "0.0502"
$ws.Range("E50").Value = "  -0.59%  "
$ws.Range("E51").Value = "  -0.23%  "

# Restore the original (default/general) style on column D so formatting
# matches the untouched cells.
$priceRange.Style = "Normal"

